$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "fundamental_data": column B (company_id) for the last two rows
# was mistakenly populated with the company-name string again; replace
# with the correct numeric ids (1 and 2).
# ---------------------------------------------------------------------
$wsFund = $wb.Worksheets.Item("fundamental_data")
$wsFund.Range("B28").Value = 1
$wsFund.Range("B29").Value = 2

# ---------------------------------------------------------------------
# Sheet "target_data": bulk-reset the base-year GHG figures (columns
# M, N, O -> base_year_ghg_s1/s2/s3) from the placeholder 200 to 0 for
# every data row (2-29).
# ---------------------------------------------------------------------
$wsTarget = $wb.Worksheets.Item("target_data")
$wsTarget.Range("M2:O27").Value = 0

# Rows 28 and 29 also need column B fixed (same issue as fundamental_data),
# the scope/metric columns (C, E) realigned to the normal "Absolute" /
# "S1+S2" choices used by every other row, the stray H column cleared, and
# the base-year GHG values reset to 0 like the rest of the table.
$wsTarget.Range("B28").Value = 1
$wsTarget.Range("B29").Value = 2

$wsTarget.Range("C2").Copy()
$wsTarget.Range("C28").PasteSpecial(-4122)
$wsTarget.Range("C29").PasteSpecial(-4122)

$wsTarget.Range("E2").Copy()
$wsTarget.Range("E28").PasteSpecial(-4122)
$wsTarget.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTarget.Range("C28").Value = "Absolute"
$wsTarget.Range("C29").Value = "Absolute"
$wsTarget.Range("E28").Value = "S1+S2"
$wsTarget.Range("E29").Value = "S1+S2"

$wsTarget.Range("H28").ClearContents()
$wsTarget.Range("H29").ClearContents()

$wsTarget.Range("M28:O29").Value = 0
$wsTarget.Range("M2").Copy()
$wsTarget.Range("M28:O29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Restore the on-screen selections that Excel recorded for each sheet.
# ---------------------------------------------------------------------
[void]$wsFund.Activate()
[void]$wsFund.Range("C28").Select()

[void]$wsTarget.Activate()
[void]$wsTarget.Range("L3").Select()
